$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Summary"
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B2").Value = 1402.79
$wsSummary.Range("E2").Value = 8597.2099999999991
$wsSummary.Range("F2").Value = 275.67
$wsSummary.Range("G2").ClearContents()

$wsSummary.Range("A5").Value = 0.23
$wsSummary.Range("B5").Value = 0.23

[void]$wsSummary.Range("B5").Select()

# ---------------------------------------------------------------------
# Sheet "Repayment schedule"
# ---------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

$wsSchedule.Range("J5").Value = 0.23
$wsSchedule.Range("K5").Value = 887.95
$wsSchedule.Range("P5").Value = 275.67

[void]$wsSchedule.Range("E12").Select()

# ---------------------------------------------------------------------
# Sheet "Transactions"
# ---------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")

# Simple id renumbering (no formatting changes involved)
$wsTrans.Range("A2").Value = 3400
$wsTrans.Range("A3").Value = 3399
$wsTrans.Range("A6").Value = 3398
$wsTrans.Range("A7").Value = 3394

# Row 3 updates
$wsTrans.Range("E3").Value = 23.24
$wsTrans.Range("I3").Value = 0.23

# Row 4 / Row 5 swap their "disbursement" vs "repayment" figures.
# F4 <-> F5 need to swap their number formats (s7 <-> s10), so stash one
# of them in a scratch cell first to avoid clobbering the source format.
$wsTrans.Range("F5").Copy() | Out-Null
$wsTrans.Range("Z1").PasteSpecial(-4122) | Out-Null   # stash F5's format (s10)

$wsTrans.Range("F4").Copy() | Out-Null
$wsTrans.Range("F5").PasteSpecial(-4122) | Out-Null   # F5 <- F4's format (s7)

$wsTrans.Range("Z1").Copy() | Out-Null
$wsTrans.Range("F4").PasteSpecial(-4122) | Out-Null   # F4 <- stashed format (s10)

# J4 needs G4's current format (s7); J5 needs E5's current format (s9).
# Neither G4 nor E5 change format, so these can be copied directly.
$wsTrans.Range("G4").Copy() | Out-Null
$wsTrans.Range("J4").PasteSpecial(-4122) | Out-Null   # J4 <- s7

$wsTrans.Range("E5").Copy() | Out-Null
$wsTrans.Range("J5").PasteSpecial(-4122) | Out-Null   # J5 <- s9

$excel.CutCopyMode = 0
$wsTrans.Range("Z1").Clear()

$wsTrans.Range("A4").Value = 3397
$wsTrans.Range("D4").Value = "Repayment"
$wsTrans.Range("E4").Value = 1500
$wsTrans.Range("F4").Value = 1402.79
$wsTrans.Range("G4").Value = 96.98
$wsTrans.Range("I4").Value = 0.23
$wsTrans.Range("J4").Value = 0

$wsTrans.Range("A5").Value = 3396
$wsTrans.Range("D5").Value = "Disbursement"
$wsTrans.Range("E5").Value = 5000
$wsTrans.Range("F5").Value = 0
$wsTrans.Range("G5").Value = 0
$wsTrans.Range("J5").Value = 10000

[void]$wsTrans.Range("D6").Select()
